$p = $ppt.ActivePresentation
$src = $p.Slides.Item(16)
$dst = $p.Slides.Item(21)

# The source deck already contains an oval callout (noFill + red outline +
# the standard "Intense Effect - Accent 1" shape style) earlier in the
# deck. Copy/paste it onto the target slide so the new shape picks up the
# exact same <p:style> theme refs and <p:txBody> placeholder paragraph
# that PowerPoint stamps on a freshly inserted AutoShape, then restyle it
# to match the new oval's own geometry/line weight.
$oval = $null
for ($i = 1; $i -le $src.Shapes.Count; $i++) {
    $candidate = $src.Shapes.Item($i)
    if ($candidate.Name -like "Oval*") {
        $oval = $candidate
        break
    }
}
if ($oval -eq $null) {
    $oval = $src.Shapes.Item(4)
}

$oval.Copy()
$pasted = $dst.Shapes.Paste()
$shp = $pasted.Item(1)
$shp.Name = "Oval 3"

# Position/size (EMU 5975797,5035640 / 2439765x373487 expressed in points).
$shp.Left = 470.53521728515625
$shp.Top = 396.5071105957031
$shp.Width = 192.10748291015625
$shp.Height = 29.40842628479004

# Outline: no fill, solid red 57150 EMU (4.5pt) line.
$shp.Fill.Visible = $false
$shp.Line.Visible = $true
$shp.Line.ForeColor.RGB = 255
$shp.Line.Weight = 4.5

# Entrance animation: Appear, on click, targeting just this shape.
$seq = $dst.TimeLine.MainSequence
$eff = $seq.AddEffect($shp, 1, 0, 1)
